$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (class 1, row id 1): replace the four course-code values with new ones.
$ws.Range("F2").Value = "K0meigYSUjKJVAITHJ7a"
$ws.Range("C2").Value = "wLCEQ5VMXX8G0AD1beYm"
$ws.Range("D2").Value = "qgp4IRO3NDUsyc2cVjR8"
$ws.Range("E2").Value = "4UkJm53t4lerabUXM7Q4"

# Row 8 (class 2, row id 7): replace the four course-code values with new ones.
$ws.Range("C8").Value = "61BMrbUIZeEfecDBJvjZ"
$ws.Range("D8").Value = "KbNuHf5ikcHQRJs9iZXr"
$ws.Range("E8").Value = "nqSDSqx2CMlja9fzamA8"
$ws.Range("F8").Value = "urB0Je8OneUBNl0zl0Bf"

# Move the active selection from F11 to F8.
$ws.Range("F8").Select()
